# "add improvements base on the readme description"
# Turn the small character/dialogue/scene example sheet into a generic
# "script command" table (kind, character, text, expression, position,
# options, image, animation) with dialogue / menu / scene(imageScene) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "start"
$ws.Name = "start"

# Pre-stamp the whole A1:H5 block with the same format already used by A1
# (reusing its existing style index instead of synthesizing new font/border/
# alignment combinations cell by cell).
$ws.Range("A1").Copy()
$ws.Range("A1:H5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row
$ws.Cells.Item(1, 1).Value = "kind"
$ws.Cells.Item(1, 2).Value = "character"
$ws.Cells.Item(1, 3).Value = "text"
$ws.Cells.Item(1, 4).Value = "expression"
$ws.Cells.Item(1, 5).Value = "position"
$ws.Cells.Item(1, 6).Value = "options"
$ws.Cells.Item(1, 7).Value = "image"
$ws.Cells.Item(1, 8).Value = "animation"

# Row 2: first dialogue line
$ws.Cells.Item(2, 1).Value = "dialogue"
$ws.Cells.Item(2, 2).Value = "John"
$ws.Cells.Item(2, 3).Value = "Hello"
$ws.Cells.Item(2, 4).Value = "happy"
$ws.Cells.Item(2, 5).Value = "left"

# Row 3: second dialogue line
$ws.Cells.Item(3, 1).Value = "dialogue"
$ws.Cells.Item(3, 2).Value = "John"
$ws.Cells.Item(3, 3).Value = "How are you?"
$ws.Cells.Item(3, 4).Value = "happy"
$ws.Cells.Item(3, 5).Value = "left"

# Row 4: menu with pipe/semicolon encoded options
$ws.Cells.Item(4, 1).Value = "menu"
$ws.Cells.Item(4, 6).Value = "option1|otherLabel;option2;option3"

# Row 5: scene change that swaps in a background image
$ws.Cells.Item(5, 1).Value = "scene"
$ws.Cells.Item(5, 7).Value = "imageScene"

# Row heights (rows 1-3 already carried ht 18.75 from the template sheet;
# give the two new rows the same height).
$ws.Rows.Item(4).RowHeight = 18.75
$ws.Rows.Item(5).RowHeight = 18.75

# Column widths: keep the original "bestFit" character width on most
# columns, but widen column F so the long options string fits.
$ws.Columns.Item(1).ColumnWidth = 12.666666666666666
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws.Columns.Item(3).ColumnWidth = 12.666666666666666
$ws.Columns.Item(4).ColumnWidth = 12.666666666666666
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 31.333333333333332
$ws.Columns.Item(7).ColumnWidth = 12.666666666666666
$ws.Columns.Item(8).ColumnWidth = 12.666666666666666
